# Rename sheet from "Refactorings 3.x to 4.x" to "Refactorings 4.0 to 4.1"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Refactorings 4.0 to 4.1"

# Move the active selection to C22 (matches <selection activeCell="C22" sqref="C22"/>)
$ws.Range("C22").Select()
